$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "student_id"/"staff_id" values for column I, rows 2-22 (one per student row)
$ids = @(
    "18-0138",
    "16-0132",
    "16-0184",
    "17-0047",
    "18-0144",
    "16-0119",
    "16-0102",
    "16-0092",
    "16-0100",
    "16-0157",
    "16-0169",
    "16-0131",
    "16-0195",
    "16-0127",
    "18-0225",
    "16-0094",
    "16-0140",
    "16-0107",
    "16-0111",
    "16-0147",
    "16-0113"
)

for ($i = 0; $i -lt $ids.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $ids[$i]
    # Match the row's existing formatting (same style as column B of that row)
    $ws.Cells.Item($row, 2).Copy()
    $ws.Cells.Item($row, 9).PasteSpecial(-4122)
}

$excel.CutCopyMode = $false

# Update the view: zoom out to 80%, scroll back to the left, and select the
# newly populated column.
$excel.ActiveWindow.Zoom = 80
$ws.Range("A1").Select()
$ws.Range("I2:I22").Select()
